$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 1158058.2
$ws_ALC.Range("I17").Value = 1999
$ws_ALC.Range("J17").Value = 1424841.1
$ws_ALC.Range("K17").Value = 5997
$ws_ALC.Range("L17").Value = 4274523.300000001
$ws_ALC.Range("M17").Value = -5829
$ws_ALC.Range("N17").Value = -4274859.300000001
$ws_ALC.Range("H49").Value = 98
$ws_ALC.Range("I49").Value = 98
$ws_ALC.Range("J49").Value = 0
$ws_ALC.Range("K49").Value = 294
$ws_ALC.Range("L49").Value = 0
$ws_ALC.Range("M49").Value = -158
$ws_ALC.Range("N49").Value = ""
$ws_ALC.Range("H76").Value = 4159
$ws_ALC.Range("I76").Value = 3323.75
$ws_ALC.Range("K76").Value = 3323.75
$ws_ALC.Range("M76").Value = -3008.75
$ws_ALC.Range("H79").Value = 4159
$ws_ALC.Range("I79").Value = 3323.75
$ws_ALC.Range("K79").Value = 3323.75
$ws_ALC.Range("M79").Value = -2231.75
$ws_ALC.Range("H80").Value = 971.88
$ws_ALC.Range("I80").Value = 807.0769
$ws_ALC.Range("J80").Value = 1150.4166
$ws_ALC.Range("K80").Value = 2421.2307
$ws_ALC.Range("L80").Value = 3451.2498
$ws_ALC.Range("M80").Value = -1423.2307
$ws_ALC.Range("N80").Value = -5447.2498
$ws_ALC.Range("H83").Value = 971.88
$ws_ALC.Range("I83").Value = 807.0769
$ws_ALC.Range("J83").Value = 1150.4166
$ws_ALC.Range("K83").Value = 7263.6921
$ws_ALC.Range("L83").Value = 10353.7494
$ws_ALC.Range("M83").Value = -2271.6921
$ws_ALC.Range("N83").Value = -20337.7494
$ws_ALC.Range("H96").Value = 3442.2727
$ws_ALC.Range("I96").Value = 1907.6666
$ws_ALC.Range("J96").Value = 10348
$ws_ALC.Range("K96").Value = 5722.9998
$ws_ALC.Range("L96").Value = 31044
$ws_ALC.Range("M96").Value = -4349.9998
$ws_ALC.Range("N96").Value = -33790
$ws_ALC.Range("H107").Value = 7499.75
$ws_ALC.Range("I107").Value = 4999
$ws_ALC.Range("J107").Value = 8333.333000000001
$ws_ALC.Range("K107").Value = 4999
$ws_ALC.Range("L107").Value = 8333.333000000001
$ws_ALC.Range("M107").Value = -3079
$ws_ALC.Range("N107").Value = -12173.333
$ws_ALC.Range("H110").Value = 49000
$ws_ALC.Range("J110").Value = 49000
$ws_ALC.Range("L110").Value = 49000
$ws_ALC.Range("N110").Value = -57180
$ws_ALC.Range("H111").Value = 19960.715
$ws_ALC.Range("I111").Value = 25845.2
$ws_ALC.Range("J111").Value = 5249.5
$ws_ALC.Range("K111").Value = 77535.60000000001
$ws_ALC.Range("L111").Value = 15748.5
$ws_ALC.Range("M111").Value = -74468.60000000001
$ws_ALC.Range("N111").Value = -21882.5
$ws_ALC.Range("H132").Value = 46881.5
$ws_ALC.Range("I132").Value = 29080.584
$ws_ALC.Range("K132").Value = 87241.75199999999
$ws_ALC.Range("M132").Value = -84711.75199999999
$ws_ALC.Range("H137").Value = 876093
$ws_ALC.Range("I137").Value = 22969.52
$ws_ALC.Range("J137").Value = 3542104
$ws_ALC.Range("K137").Value = 68908.56
$ws_ALC.Range("L137").Value = 10626312
$ws_ALC.Range("M137").Value = -66358.56
$ws_ALC.Range("N137").Value = -10631412

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 2131.5
$ws_ARM.Range("I2").Value = 2164.375
$ws_ARM.Range("K2").Value = 2164.375
$ws_ARM.Range("M2").Value = -2051.375
$ws_ARM.Range("H32").Value = 3972.4856
$ws_ARM.Range("I32").Value = 3663.6897
$ws_ARM.Range("K32").Value = 3663.6897
$ws_ARM.Range("M32").Value = -3376.6897
$ws_ARM.Range("H45").Value = 21224.312
$ws_ARM.Range("I45").Value = 18171.23
$ws_ARM.Range("J45").Value = 34454.332
$ws_ARM.Range("K45").Value = 18171.23
$ws_ARM.Range("L45").Value = 34454.332
$ws_ARM.Range("M45").Value = -17794.23
$ws_ARM.Range("N45").Value = -35208.332
$ws_ARM.Range("H88").Value = 3116.375
$ws_ARM.Range("I88").Value = 3125.8
$ws_ARM.Range("J88").Value = 3100.6667
$ws_ARM.Range("K88").Value = 3125.8
$ws_ARM.Range("L88").Value = 3100.6667
$ws_ARM.Range("M88").Value = -2719.8
$ws_ARM.Range("N88").Value = -3912.6667
$ws_ARM.Range("H91").Value = 3116.375
$ws_ARM.Range("I91").Value = 3125.8
$ws_ARM.Range("J91").Value = 3100.6667
$ws_ARM.Range("K91").Value = 3125.8
$ws_ARM.Range("L91").Value = 3100.6667
$ws_ARM.Range("M91").Value = -1721.8
$ws_ARM.Range("N91").Value = -5908.6667
$ws_ARM.Range("H110").Value = 11276.667
$ws_ARM.Range("I110").Value = 11587.059
$ws_ARM.Range("K110").Value = 11587.059
$ws_ARM.Range("M110").Value = -9542.058999999999
$ws_ARM.Range("H116").Value = 2131.5
$ws_ARM.Range("I116").Value = 2164.375
$ws_ARM.Range("K116").Value = 2164.375
$ws_ARM.Range("M116").Value = 129.625
$ws_ARM.Range("H122").Value = 2245.9092
$ws_ARM.Range("I122").Value = 2270.6
$ws_ARM.Range("K122").Value = 6811.799999999999
$ws_ARM.Range("M122").Value = -4361.799999999999
$ws_ARM.Range("H132").Value = 3147.6128
$ws_ARM.Range("I132").Value = 2671.0435
$ws_ARM.Range("K132").Value = 8013.130500000001
$ws_ARM.Range("M132").Value = -5483.130500000001

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 2131.5
$ws_BSM.Range("I3").Value = 2164.375
$ws_BSM.Range("K3").Value = 2164.375
$ws_BSM.Range("M3").Value = -2050.375
$ws_BSM.Range("H22").Value = 318.75
$ws_BSM.Range("I22").Value = 348
$ws_BSM.Range("J22").Value = 289.5
$ws_BSM.Range("K22").Value = 348
$ws_BSM.Range("L22").Value = 289.5
$ws_BSM.Range("M22").Value = -175
$ws_BSM.Range("N22").Value = -635.5
$ws_BSM.Range("H54").Value = 16833
$ws_BSM.Range("I54").Value = 16833
$ws_BSM.Range("K54").Value = 16833
$ws_BSM.Range("M54").Value = -16349
$ws_BSM.Range("H86").Value = 3464
$ws_BSM.Range("I86").Value = 2217.238
$ws_BSM.Range("K86").Value = 2217.238
$ws_BSM.Range("M86").Value = -1094.238
$ws_BSM.Range("H89").Value = 3464
$ws_BSM.Range("I89").Value = 2217.238
$ws_BSM.Range("K89").Value = 11086.19
$ws_BSM.Range("M89").Value = -5470.189999999999
$ws_BSM.Range("H94").Value = 1645.7333
$ws_BSM.Range("I94").Value = 1851.5
$ws_BSM.Range("J94").Value = 1234.2
$ws_BSM.Range("K94").Value = 1851.5
$ws_BSM.Range("L94").Value = 1234.2
$ws_BSM.Range("M94").Value = -1400.5
$ws_BSM.Range("N94").Value = -2136.2
$ws_BSM.Range("H99").Value = 6055.9375
$ws_BSM.Range("I99").Value = 2870
$ws_BSM.Range("J99").Value = 10152.143
$ws_BSM.Range("K99").Value = 2870
$ws_BSM.Range("L99").Value = 10152.143
$ws_BSM.Range("M99").Value = -1372
$ws_BSM.Range("N99").Value = -13148.143
$ws_BSM.Range("H107").Value = 2859.5
$ws_BSM.Range("I107").Value = 2493.3333
$ws_BSM.Range("J107").Value = 3408.75
$ws_BSM.Range("K107").Value = 2493.3333
$ws_BSM.Range("L107").Value = 3408.75
$ws_BSM.Range("M107").Value = -573.3332999999998
$ws_BSM.Range("N107").Value = -7248.75
$ws_BSM.Range("H134").Value = 1750.5834
$ws_BSM.Range("I134").Value = 1562.3823
$ws_BSM.Range("K134").Value = 4687.1469
$ws_BSM.Range("M134").Value = -2152.1469

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H9").Value = 521500
$ws_CRP.Range("J9").Value = 521500
$ws_CRP.Range("L9").Value = 521500
$ws_CRP.Range("N9").Value = -521836
$ws_CRP.Range("H31").Value = 2816934.5
$ws_CRP.Range("I31").Value = 4086101.5
$ws_CRP.Range("K31").Value = 4086101.5
$ws_CRP.Range("M31").Value = -4085806.5
$ws_CRP.Range("H34").Value = 2816934.5
$ws_CRP.Range("I34").Value = 4086101.5
$ws_CRP.Range("K34").Value = 4086101.5
$ws_CRP.Range("M34").Value = -4085899.5
$ws_CRP.Range("H35").Value = 918.75
$ws_CRP.Range("I35").Value = 925
$ws_CRP.Range("K35").Value = 925
$ws_CRP.Range("M35").Value = -631
$ws_CRP.Range("H58").Value = 4164.522
$ws_CRP.Range("I58").Value = 3433.7856
$ws_CRP.Range("K58").Value = 3433.7856
$ws_CRP.Range("M58").Value = -3230.7856
$ws_CRP.Range("H99").Value = 4959.3335
$ws_CRP.Range("I99").Value = 4959.3335
$ws_CRP.Range("J99").Value = 0
$ws_CRP.Range("K99").Value = 4959.3335
$ws_CRP.Range("L99").Value = 0
$ws_CRP.Range("M99").Value = -3461.3335
$ws_CRP.Range("N99").Value = ""
$ws_CRP.Range("H109").Value = 35450
$ws_CRP.Range("J109").Value = 35450
$ws_CRP.Range("L109").Value = 35450
$ws_CRP.Range("N109").Value = -37530
$ws_CRP.Range("H123").Value = 0
$ws_CRP.Range("J123").Value = 0
$ws_CRP.Range("L123").Value = 0
$ws_CRP.Range("N123").Value = ""
$ws_CRP.Range("H126").Value = 4959.3335
$ws_CRP.Range("I126").Value = 4959.3335
$ws_CRP.Range("J126").Value = 0
$ws_CRP.Range("K126").Value = 14878.0005
$ws_CRP.Range("L126").Value = 0
$ws_CRP.Range("M126").Value = -12408.0005
$ws_CRP.Range("N126").Value = ""
$ws_CRP.Range("H134").Value = 21298.969
$ws_CRP.Range("I134").Value = 24474.889
$ws_CRP.Range("J134").Value = 4149
$ws_CRP.Range("K134").Value = 73424.667
$ws_CRP.Range("L134").Value = 12447
$ws_CRP.Range("M134").Value = -70889.667
$ws_CRP.Range("N134").Value = -17517
$ws_CRP.Range("H136").Value = 4164.522
$ws_CRP.Range("I136").Value = 3433.7856
$ws_CRP.Range("K136").Value = 10301.3568
$ws_CRP.Range("M136").Value = -7751.356800000001

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H11").Value = 83428.25
$ws_CUL.Range("I11").Value = 85.71429000000001
$ws_CUL.Range("J11").Value = 200107.8
$ws_CUL.Range("K11").Value = 257.14287
$ws_CUL.Range("L11").Value = 600323.3999999999
$ws_CUL.Range("M11").Value = -117.14287
$ws_CUL.Range("N11").Value = -600603.3999999999
$ws_CUL.Range("H40").Value = 120.333336
$ws_CUL.Range("I40").Value = 129.125
$ws_CUL.Range("K40").Value = 516.5
$ws_CUL.Range("M40").Value = -447.5
$ws_CUL.Range("H88").Value = 66670000
$ws_CUL.Range("J88").Value = 100000000
$ws_CUL.Range("L88").Value = 300000000
$ws_CUL.Range("N88").Value = -300000856
$ws_CUL.Range("H91").Value = 66670000
$ws_CUL.Range("J91").Value = 100000000
$ws_CUL.Range("L91").Value = 300000000
$ws_CUL.Range("N91").Value = -300002964
$ws_CUL.Range("H96").Value = 10250
$ws_CUL.Range("H107").Value = 33334200
$ws_CUL.Range("I107").Value = 41667176
$ws_CUL.Range("J107").Value = 2300
$ws_CUL.Range("K107").Value = 125001528
$ws_CUL.Range("L107").Value = 6900
$ws_CUL.Range("M107").Value = -124999608
$ws_CUL.Range("N107").Value = -10740
$ws_CUL.Range("H117").Value = 989.44446
$ws_CUL.Range("I117").Value = 557.8570999999999
$ws_CUL.Range("J117").Value = 2500
$ws_CUL.Range("K117").Value = 1673.5713
$ws_CUL.Range("L117").Value = 7500
$ws_CUL.Range("M117").Value = 1768.4287
$ws_CUL.Range("N117").Value = -14384
$ws_CUL.Range("H121").Value = 2952.3845
$ws_CUL.Range("I121").Value = 2134.5264
$ws_CUL.Range("J121").Value = 5172.2856
$ws_CUL.Range("K121").Value = 6403.5792
$ws_CUL.Range("L121").Value = 15516.8568
$ws_CUL.Range("M121").Value = -5093.5792
$ws_CUL.Range("N121").Value = -18136.8568
$ws_CUL.Range("H131").Value = 18567.6
$ws_CUL.Range("J131").Value = 2066.0784
$ws_CUL.Range("L131").Value = 6198.235199999999
$ws_CUL.Range("N131").Value = -16278.2352

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 5841.778
$ws_GSM.Range("I70").Value = 5466.6665
$ws_GSM.Range("J70").Value = 6029.3335
$ws_GSM.Range("K70").Value = 5466.6665
$ws_GSM.Range("L70").Value = 6029.3335
$ws_GSM.Range("M70").Value = -5196.6665
$ws_GSM.Range("N70").Value = -6569.3335
$ws_GSM.Range("H73").Value = 5841.778
$ws_GSM.Range("I73").Value = 5466.6665
$ws_GSM.Range("J73").Value = 6029.3335
$ws_GSM.Range("K73").Value = 5466.6665
$ws_GSM.Range("L73").Value = 6029.3335
$ws_GSM.Range("M73").Value = -4530.6665
$ws_GSM.Range("N73").Value = -7901.3335
$ws_GSM.Range("H122").Value = 6043.926
$ws_GSM.Range("I122").Value = 6043.926
$ws_GSM.Range("K122").Value = 18131.778
$ws_GSM.Range("M122").Value = -15681.778
$ws_GSM.Range("H126").Value = 22498.766
$ws_GSM.Range("I126").Value = 42700
$ws_GSM.Range("J126").Value = 4542.1113
$ws_GSM.Range("K126").Value = 128100
$ws_GSM.Range("L126").Value = 13626.3339
$ws_GSM.Range("M126").Value = -125630
$ws_GSM.Range("N126").Value = -18566.3339
$ws_GSM.Range("H132").Value = 17837.906
$ws_GSM.Range("I132").Value = 18760.434
$ws_GSM.Range("K132").Value = 56281.302
$ws_GSM.Range("M132").Value = -53751.302
$ws_GSM.Range("H136").Value = 42445.2
$ws_GSM.Range("J136").Value = 42445.2
$ws_GSM.Range("L136").Value = 127335.6
$ws_GSM.Range("N136").Value = -132435.6

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 11857.917
$ws_LTW.Range("I7").Value = 12779.5
$ws_LTW.Range("J7").Value = 7250
$ws_LTW.Range("K7").Value = 12779.5
$ws_LTW.Range("L7").Value = 7250
$ws_LTW.Range("M7").Value = -12667.5
$ws_LTW.Range("N7").Value = -7474
$ws_LTW.Range("H16").Value = 3977.0715
$ws_LTW.Range("I16").Value = 3699.0908
$ws_LTW.Range("K16").Value = 3699.0908
$ws_LTW.Range("M16").Value = -3529.0908
$ws_LTW.Range("H22").Value = 778.8889
$ws_LTW.Range("I22").Value = 778.8889
$ws_LTW.Range("K22").Value = 778.8889
$ws_LTW.Range("M22").Value = -483.8889
$ws_LTW.Range("H27").Value = 778.8889
$ws_LTW.Range("I27").Value = 778.8889
$ws_LTW.Range("K27").Value = 778.8889
$ws_LTW.Range("M27").Value = -671.8889
$ws_LTW.Range("H40").Value = 7565.9165
$ws_LTW.Range("I40").Value = 5799.3335
$ws_LTW.Range("K40").Value = 5799.3335
$ws_LTW.Range("M40").Value = -5663.3335
$ws_LTW.Range("H100").Value = 3247.6667
$ws_LTW.Range("I100").Value = 2996.3333
$ws_LTW.Range("K100").Value = 2996.3333
$ws_LTW.Range("M100").Value = -2455.3333
$ws_LTW.Range("H126").Value = 11857.917
$ws_LTW.Range("I126").Value = 12779.5
$ws_LTW.Range("J126").Value = 7250
$ws_LTW.Range("K126").Value = 38338.5
$ws_LTW.Range("L126").Value = 21750
$ws_LTW.Range("M126").Value = -35868.5
$ws_LTW.Range("N126").Value = -26690
$ws_LTW.Range("H132").Value = 5147.923
$ws_LTW.Range("I132").Value = 4321.143
$ws_LTW.Range("K132").Value = 12963.429
$ws_LTW.Range("M132").Value = -10433.429
$ws_LTW.Range("H136").Value = 3012.3333
$ws_LTW.Range("I136").Value = 3104.3635
$ws_LTW.Range("J136").Value = 2000
$ws_LTW.Range("K136").Value = 9313.0905
$ws_LTW.Range("L136").Value = 6000
$ws_LTW.Range("M136").Value = -6763.0905
$ws_LTW.Range("N136").Value = -11100

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H24").Value = 11666.667
$ws_WVR.Range("I24").Value = 10000
$ws_WVR.Range("K24").Value = 10000
$ws_WVR.Range("M24").Value = -9770
$ws_WVR.Range("H30").Value = 10000
$ws_WVR.Range("J30").Value = 0
$ws_WVR.Range("L30").Value = 0
$ws_WVR.Range("N30").Value = ""
$ws_WVR.Range("H49").Value = 23000
$ws_WVR.Range("J49").Value = 23000
$ws_WVR.Range("L49").Value = 23000
$ws_WVR.Range("N49").Value = -23460
$ws_WVR.Range("H92").Value = 0
$ws_WVR.Range("J92").Value = 0
$ws_WVR.Range("L92").Value = 0
$ws_WVR.Range("N92").Value = ""
$ws_WVR.Range("H122").Value = 4229.4517
$ws_WVR.Range("I122").Value = 2646.64
$ws_WVR.Range("J122").Value = 10824.5
$ws_WVR.Range("K122").Value = 7939.92
$ws_WVR.Range("L122").Value = 32473.5
$ws_WVR.Range("M122").Value = -5489.92
$ws_WVR.Range("N122").Value = -37373.5
$ws_WVR.Range("H126").Value = 3928.7368
$ws_WVR.Range("I126").Value = 2977.875
$ws_WVR.Range("K126").Value = 8933.625
$ws_WVR.Range("M126").Value = -6463.625
$ws_WVR.Range("H132").Value = 9632.695
$ws_WVR.Range("I132").Value = 7697.6665
$ws_WVR.Range("K132").Value = 23092.9995
$ws_WVR.Range("M132").Value = -20562.9995
$ws_WVR.Range("H136").Value = 1708.762
$ws_WVR.Range("I136").Value = 1346.2941
$ws_WVR.Range("K136").Value = 4038.8823
$ws_WVR.Range("M136").Value = -1488.8823
$ws_WVR.Range("H137").Value = 107500
$ws_WVR.Range("J137").Value = 107500
$ws_WVR.Range("L137").Value = 107500
$ws_WVR.Range("N137").Value = -117700
